$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values first
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 0.8557439673732903

# Delete rows 3 through 17 (the remaining data rows)
$ws.Range("A3:B17").EntireRow.Delete()
